$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.968.51'
$ws.Range("E2").Value = '  -0.03%  '

$ws.Range("D3").Value = '3.412.56'
$ws.Range("E3").Value = '  -0.13%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").Value = "'408.71"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("D6").Value = "'129.00"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -4.49%  '

$ws.Range("D7").Value = "'0.641"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +8.69%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = "'0.733"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +6.97%  '

$ws.Range("D10").Value = "'0.141"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +16.11%  '

$ws.Range("D11").Value = "'42.37"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.53%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = "'0.141"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = "'0.0000213"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +63.67%  '

$ws.Range("D14").Value = '3.959.03'
$ws.Range("E14").Value = '  +0.20%  '

$ws.Range("D15").Value = "'8.98"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +6.54%  '

$ws.Range("D16").Value = "'20.79"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +4.50%  '

$ws.Range("D17").Value = '3.375.35'
$ws.Range("E17").Value = '  -1.83%  '

$ws.Range("D18").Value = "'12.11"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +9.87%  '

$ws.Range("D19").Value = "'1.06"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +4.63%  '

$ws.Range("D20").Value = '61.914.42'
$ws.Range("E20").Value = '  -0.06%  '

$ws.Range("D21").Value = "'443.15"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +41.56%  '

$ws.Range("D22").Value = "'91.33"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +7.61%  '

$ws.Range("E23").Value = '  -0.98%  '

$ws.Range("D24").Value = "'13.01"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +1.29%  '

$ws.Range("D25").Value = "'3.23"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +2.86%  '

$ws.Range("D26").Value = "'33.88"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +14.25%  '

$ws.Range("D27").Value = "'8.75"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +7.46%  '

$ws.Range("E28").Value = '  +0.22%  '

$ws.Range("D29").Value = "'7.61"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +0.95%  '

$ws.Range("D30").Value = "'2.67"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -5.24%  '

$ws.Range("E31").Value = '  +5.69%  '

$ws.Range("E32").Value = '  -2.29%  '

$ws.Range("E33").Value = '  -0.79%  '

$ws.Range("D34").Value = "'42.44"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").Value = "'0.0502"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +4.13%  '

$ws.Range("D37").Value = "'53.71"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +3.79%  '

$ws.Range("E38").Value = '  +0.02%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").Value = "'0.135"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +8.21%  '

$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").Value = "'3.35"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -1.41%  '

$ws.Range("D41").Value = "'2.90"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -1.31%  '

$ws.Range("D42").Value = "'0.313"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +2.55%  '

$ws.Range("D43").Value = "'141.14"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +2.07%  '

$ws.Range("D44").Value = "'1.96"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -0.77%  '

$ws.Range("D45").Value = "'4.08"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +1.38%  '

$ws.Range("D46").Value = "'2.40"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +8.37%  '

$ws.Range("E47").Value = '  -1.06%  '

$ws.Range("D48").Value = "'22.33"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +4.77%  '

$ws.Range("D49").Value = '3.763.16'
$ws.Range("E49").Value = '  +0.34%  '

$ws.Range("D50").Value = '2.102.49'
$ws.Range("E50").Value = '  -0.76%  '

$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").Value = "'105.81"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +27.12%  '
